# "Committing to new branch" re-run of the breed-filter stat query:
# StatOutput / StatOutput_Message now reflect the *Gordon Setter* stat
# query (instead of the old Akita one) -> number_of_sample/cases/study
# go from 2/1/1 to 3/2/2, and the logged Cypher text on
# StatOutput_Message!A18 is updated to match.

$wb = $excel.ActiveWorkbook

$statOutput = $wb.Worksheets.Item("StatOutput")

# StatOutput row 2: A2=number_of_files(1, unchanged), B2=number_of_sample,
# C2=number_of_cases, D2=number_of_study.
# C2 and D2 both pick up the value that used to live in B2 ("2"), so copy
# it over first (Copy keeps the cell as a text/shared-string cell, same
# as the source, with no formatting side effects).
$statOutput.Range("B2").Copy($statOutput.Range("C2"))
$statOutput.Range("B2").Copy($statOutput.Range("D2"))

# B2 itself becomes "3". Stage the new text on a scratch sheet (forcing
# text via NumberFormat "@" so it isn't silently read back as a number),
# copy it into B2, then strip the borrowed number-format and remove the
# scratch sheet so nothing else on the workbook is touched.
$scratchSheet = $wb.Worksheets.Add()
$statOutput = $wb.Worksheets.Item("StatOutput")
$scratch = $scratchSheet.Range("A1")
$scratch.NumberFormat = "@"
$scratch.Value = "3"
$scratch.Copy($statOutput.Range("B2"))
$statOutput.Range("B2").ClearFormats()
$scratchSheet.Delete() | Out-Null

# StatOutput_Message!A18 holds the logged Cypher query text for the stat
# run; swap the Akita filter for the Gordon Setter one (the OPTIONAL
# MATCH file/sample-counting variant of the query).
$statMessage = $wb.Worksheets.Item("StatOutput_Message")
$gordonQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Gordon Setter']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$statMessage.Range("A18").Value = $gordonQuery
